# Edit applied:
#  1. Slide 6's table switches from the custom "Table_0" table style
#     ({BE9F6816-7DFA-4CBD-857A-54C6BFAB2BAD}) to the built-in table style
#     {282C953C-A6C3-4C51-8CF9-BFC06363782B} (PowerPoint Table Styles gallery).
#  2. The presentation's theme color scheme is switched from the "Integral"
#     palette over to the stock "Office" palette (the deck's look was
#     changed from the Integral theme to the default Office Theme).

$p = $ppt.ActivePresentation

# --- 1. Update the table style on slide 6 ---------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{282C953C-A6C3-4C51-8CF9-BFC06363782B}")
    }
}

# --- 2. Swap the active theme's 12 scheme colors to the Office palette ----
$officeColors = @{
    1  = 0x000000   # dk1
    2  = 0xFFFFFF   # lt1
    3  = 0x44546A   # dk2
    4  = 0xE7E6E6   # lt2
    5  = 0x5B9BD5   # accent1
    6  = 0xED7D31   # accent2
    7  = 0xA5A5A5   # accent3
    8  = 0xFFC000   # accent4
    9  = 0x4472C4   # accent5
    10 = 0x70AD47   # accent6
    11 = 0x0563C1   # hlink
    12 = 0x954F72   # folHlink
}

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $hex = $officeColors[$i]
    $r = [math]::Floor($hex / 0x10000) -band 0xFF
    $g = [math]::Floor($hex / 0x100) -band 0xFF
    $b = $hex -band 0xFF
    $tcs.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
